$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded; it belongs chronologically
# right after the existing row 81, so insert a fresh row at 82 and push the
# old rows 82:87 down to 83:88 (their own values are untouched by the shift).
$ws.Rows("82:82").Insert()

# Populate the newly inserted row 82 with the new observation.
$row = 82
$ws.Cells.Item($row, 1).Value = 8
$ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 44610
$ws.Cells.Item($row, 5).Value = 4
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100109
$ws.Cells.Item($row, 8).Value = "Uva"
$ws.Cells.Item($row, 9).Value = 100109001
$ws.Cells.Item($row, 10).Value = "Uva"
$ws.Cells.Item($row, 11).Value = "Red Globe"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 400
$ws.Cells.Item($row, 14).Value = 11500
$ws.Cells.Item($row, 15).Value = 12000
$ws.Cells.Item($row, 16).Value = 11750
$ws.Cells.Item($row, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item($row, 18).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 19).Value = 653
$ws.Cells.Item($row, 20).Value = 18
